# Thu, Apr 02, 2020  4:06:23 PM
#
# 1) Re-point the table on slide 6 at the built-in table style
#    {E4C819CD-591F-45A4-B32C-BF7DFBF22868} instead of the custom
#    theme-based style {9D2D8916-49AE-4024-B8AF-26594AD2762A}.
# 2) Swap the presentation's theme color palette from the "Integral"
#    palette over to the standard "Office" palette (dk1/lt1/dk2/lt2/
#    accent1-6/hlink/folHlink), which is what the deck's second theme
#    part already held.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$tableShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTable) {
            $tableShape = $shp
        }
    }
}

if ($tableShape -ne $null) {
    $tableShape.Table.ApplyStyle("{E4C819CD-591F-45A4-B32C-BF7DFBF22868}")
}

# --- 2. Theme colors --------------------------------------------------
# Index order for ThemeColorScheme.Item(): dk1, lt1, dk2, lt2,
# accent1..accent6, hlink, folHlink (MsoThemeColorSchemeIndex order).
# RGB is packed the VBA way: r + g*256 + b*65536.
$officeRgb = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme
for ($i = 1; $i -le $scheme.Count; $i++) {
    $scheme.Item($i).RGB = $officeRgb[$i - 1]
}
